# Applies the commit's worksheet changes:
#  - test_b (sheet1): add columns D:G replicating column C's data pattern
#    with headers S2/S3/S4/S5, plus hyperlinks on row 3 (D3:G3) pointing at
#    the same w3school demo URL as C3.
#  - test_b becomes the active (selected) sheet/tab instead of test_a.
#  - test_a (sheet2) is no longer the active tab.

$wb = $excel.ActiveWorkbook
$wsB = $wb.Worksheets.Item("test_b")

$url = "https://www.w3school.com.cn/tiy/t.asp?f=jquery_fadeout"

# New header labels for the added columns.
$headers = @("S2", "S3", "S4", "S5")
$cols = @("D", "E", "F", "G")

for ($i = 0; $i -lt 4; $i++) {
    $col = $cols[$i]

    # Header row (row 1) - same style as the rest of row 1.
    $wsB.Range($col + "1").Value = $headers[$i]

    # Rows 2, 4, 5 simply repeat the value already in column C.
    $wsB.Range($col + "2").Value = $wsB.Range("C2").Value()
    $wsB.Range($col + "4").Value = $wsB.Range("C4").Value()
    $wsB.Range($col + "5").Value = $wsB.Range("C5").Value()

    # Row 3 is the hyperlink cell - set value + hyperlink, then restore the
    # same number format / hyperlink style that column C3 already has.
    $wsB.Range($col + "3").Value = $url
    $wsB.Hyperlinks.Add($wsB.Range($col + "3"), $url)
    $wsB.Range("C3").Copy()
    $wsB.Range($col + "3").PasteSpecial(-4122)
}

# test_b is now the selected/active sheet (tabSelected moves from test_a).
$wsB.Activate()
$wsB.Range("G12").Select()
